$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: add P1=14, Q1=15 (continuing the existing 0..13 sequence).
# Copy formatting from the last existing header cell (O1) so the new cells
# get the same style (bold, centered, bordered) as the rest of row 1.
$ws.Cells.Item(1, 15).Copy($ws.Cells.Item(1, 16))
$ws.Cells.Item(1, 16).Value = 14   # P1

$ws.Cells.Item(1, 15).Copy($ws.Cells.Item(1, 17))
$ws.Cells.Item(1, 17).Value = 15   # Q1

# Rows 2..25: swap values in columns I/K and M/O, and add new columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2
}
